# Natmi following Dr Hou advice
# Update Agrn-Musk LR-pairs sheet: recompute rows 2-7 and add new rows 8-10
# covering the full 3x3 cross-product of clusters (ECs, FAPs, sCs) as
# sending/target pairs (previously only off-diagonal pairs were present).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Agrn"
$ws.Range("C2").Value = "Musk"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 17.306265
$ws.Range("H2").Value = 51.918795
$ws.Range("I2").Value = 0.5463168539988408
$ws.Range("J2").Value = 0.5463168539988407
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1421586666666667
$ws.Range("N2").Value = 0.426476
$ws.Range("O2").Value = 0.008210963676175792
$ws.Range("P2").Value = 0.008210963676175792
$ws.Range("Q2").Value = 2.46023555738
$ws.Range("R2").Value = 22.14212001642
$ws.Range("S2").Value = 0.004485787843867116
$ws.Range("T2").Value = 0.004485787843867115

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Agrn"
$ws.Range("C3").Value = "Musk"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 17.306265
$ws.Range("H3").Value = 51.918795
$ws.Range("I3").Value = 0.5463168539988408
$ws.Range("J3").Value = 0.5463168539988407
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.072364333333334
$ws.Range("N3").Value = 18.217093
$ws.Range("O3").Value = 0.3507345991533317
$ws.Range("P3").Value = 0.3507345991533317
$ws.Range("Q3").Value = 105.089946329215
$ws.Range("R3").Value = 945.8095169629352
$ws.Range("S3").Value = 0.1916122227979927
$ws.Range("T3").Value = 0.1916122227979926

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Agrn"
$ws.Range("C4").Value = "Musk"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 17.306265
$ws.Range("H4").Value = 51.918795
$ws.Range("I4").Value = 0.5463168539988408
$ws.Range("J4").Value = 0.5463168539988407
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 11.09875133333333
$ws.Range("N4").Value = 33.296254
$ws.Range("O4").Value = 0.6410544371704925
$ws.Range("P4").Value = 0.6410544371704925
$ws.Range("Q4").Value = 192.07793174377
$ws.Range("R4").Value = 1728.70138569393
$ws.Range("S4").Value = 0.350218843356981
$ws.Range("T4").Value = 0.3502188433569809

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Agrn"
$ws.Range("C5").Value = "Musk"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.344413333333334
$ws.Range("H5").Value = 10.03324
$ws.Range("I5").Value = 0.1055750256186672
$ws.Range("J5").Value = 0.1055750256186672
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1421586666666667
$ws.Range("N5").Value = 0.426476
$ws.Range("O5").Value = 0.008210963676175792
$ws.Range("P5").Value = 0.008210963676175792
$ws.Range("Q5").Value = 0.475437340248889
$ws.Range("R5").Value = 4.278936062240001
$ws.Range("S5").Value = 0.0008668727004662049
$ws.Range("T5").Value = 0.0008668727004662048

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Agrn"
$ws.Range("C6").Value = "Musk"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.344413333333334
$ws.Range("H6").Value = 10.03324
$ws.Range("I6").Value = 0.1055750256186672
$ws.Range("J6").Value = 0.1055750256186672
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.072364333333334
$ws.Range("N6").Value = 18.217093
$ws.Range("O6").Value = 0.3507345991533317
$ws.Range("P6").Value = 0.3507345991533317
$ws.Range("Q6").Value = 20.30849624125778
$ws.Range("R6").Value = 182.77646617132
$ws.Range("S6").Value = 0.03702881429096596
$ws.Range("T6").Value = 0.03702881429096595

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Agrn"
$ws.Range("C7").Value = "Musk"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.344413333333334
$ws.Range("H7").Value = 10.03324
$ws.Range("I7").Value = 0.1055750256186672
$ws.Range("J7").Value = 0.1055750256186672
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 11.09875133333333
$ws.Range("N7").Value = 33.296254
$ws.Range("O7").Value = 0.6410544371704925
$ws.Range("P7").Value = 0.6410544371704925
$ws.Range("Q7").Value = 37.11881194255112
$ws.Range("R7").Value = 334.0693074829601
$ws.Range("S7").Value = 0.067679338627235
$ws.Range("T7").Value = 0.067679338627235

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Agrn"
$ws.Range("C8").Value = "Musk"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 11.02739433333333
$ws.Range("H8").Value = 33.082183
$ws.Range("I8").Value = 0.3481081203824922
$ws.Range("J8").Value = 0.3481081203824921
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1421586666666667
$ws.Range("N8").Value = 0.426476
$ws.Range("O8").Value = 0.008210963676175792
$ws.Range("P8").Value = 0.008210963676175792
$ws.Range("Q8").Value = 1.567639675234223
$ws.Range("R8").Value = 14.108757077108
$ws.Range("S8").Value = 0.002858303131842473
$ws.Range("T8").Value = 0.002858303131842473

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Agrn"
$ws.Range("C9").Value = "Musk"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 11.02739433333333
$ws.Range("H9").Value = 33.082183
$ws.Range("I9").Value = 0.3481081203824922
$ws.Range("J9").Value = 0.3481081203824921
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.072364333333334
$ws.Range("N9").Value = 18.217093
$ws.Range("O9").Value = 0.3507345991533317
$ws.Range("P9").Value = 0.3507345991533317
$ws.Range("Q9").Value = 66.96235603933545
$ws.Range("R9").Value = 602.6612043540191
$ws.Range("S9").Value = 0.1220935620643731
$ws.Range("T9").Value = 0.1220935620643731

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Agrn"
$ws.Range("C10").Value = "Musk"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 11.02739433333333
$ws.Range("H10").Value = 33.082183
$ws.Range("I10").Value = 0.3481081203824922
$ws.Range("J10").Value = 0.3481081203824921
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 11.09875133333333
$ws.Range("N10").Value = 33.296254
$ws.Range("O10").Value = 0.6410544371704925
$ws.Range("P10").Value = 0.6410544371704925
$ws.Range("Q10").Value = 122.3903075602758
$ws.Range("R10").Value = 1101.512768042482
$ws.Range("S10").Value = 0.2231562551862766
$ws.Range("T10").Value = 0.2231562551862765

